$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 406; existing rows 406-469 shift down to 407-470.
$ws.Rows.Item(406).Insert()

# Populate the newly inserted row 406 with the new data record.
$ws.Cells.Item(406, 1).Value = 6
$ws.Cells.Item(406, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(406, 3).Value = "Metropolitana"
$ws.Cells.Item(406, 4).Value = 44776
$ws.Cells.Item(406, 5).Value = 13
$ws.Cells.Item(406, 6).Value = 100112043
$ws.Cells.Item(406, 7).Value = "Pepino ensalada"
$ws.Cells.Item(406, 8).Value = "Sin especificar"
$ws.Cells.Item(406, 9).Value = "Primera"
$ws.Cells.Item(406, 10).Value = 400
$ws.Cells.Item(406, 11).Value = 17000
$ws.Cells.Item(406, 12).Value = 19000
$ws.Cells.Item(406, 13).Value = 17850
$ws.Cells.Item(406, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(406, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(406, 16).Value = 298
$ws.Cells.Item(406, 17).Value = 60
$ws.Cells.Item(406, 18).Value = "Hortaliza"

# Match the date column's number format (style) used throughout column D.
$ws.Cells.Item(406, 4).NumberFormat = $ws.Cells.Item(407, 4).NumberFormat
